# project plan update with status
# Adds a "status" column (D) to the Task Distribution table on Sheet2:
# Phase 1 rows get Done/Done/WIP markers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New status column for the Phase-1 task-distribution rows.
$ws.Range("D4").Value = "Done"
$ws.Range("D5").Value = "Done"
$ws.Range("D6").Value = "WIP"
# Row 7 (Ruturaj / Setup MLFlow for experiment tracking) is left without a status,
# matching the source edit.

# The description column (C) for every task row was carrying a redundant,
# duplicate cell format (explicit "apply default font" on top of the normal
# left/vcenter/indent alignment). Touching the font with its current (already
# unbold) value collapses each of these cells back onto the canonical shared
# style instead of the redundant duplicate one.
$ws.Range("C4:C7").Font.Bold = $false
$ws.Range("C9:C12").Font.Bold = $false
$ws.Range("C14:C17").Font.Bold = $false

# Move the active selection the way it was left after the edit.
$ws.Range("C19").Select() | Out-Null
